$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is an unambiguous plain number (e.g. "179.85") would be
# auto-converted from text to a numeric type by Excel (losing fixed decimal/
# trailing-zero formatting, e.g. "26.50" -> 26.5). Force those specific cells to
# keep a Text number format first so the values are stored as exact text, matching
# the source data (which always stores these as text strings).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = '70.890.85'
$ws.Range("E2").Value = '  +3.56%  '
$ws.Range("D3").Value = '2.611.24'
$ws.Range("E3").Value = '  +3.76%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '604.52'
$ws.Range("E5").Value = '  +2.13%  '
$ws.Range("D6").Value = '179.85'
$ws.Range("E6").Value = '  +2.24%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  +1.86%  '
$ws.Range("D9").Value = '2.611.17'
$ws.Range("E9").Value = '  +3.83%  '
$ws.Range("E10").Value = '  +13.87%  '
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("E12").Value = '  +2.93%  '
$ws.Range("D13").Value = '5.01'
$ws.Range("E13").Value = '  +0.50%  '
$ws.Range("D14").Value = '3.098.21'
$ws.Range("E14").Value = '  +1.73%  '
$ws.Range("D15").Value = '26.50'
$ws.Range("E15").Value = '  +2.82%  '
$ws.Range("E16").Value = '  +6.68%  '
$ws.Range("D17").Value = '70.907.25'
$ws.Range("E17").Value = '  +3.87%  '
$ws.Range("D18").Value = '2.637.02'
$ws.Range("E18").Value = '  +5.06%  '
$ws.Range("D19").Value = '379.39'
$ws.Range("E19").Value = '  +7.96%  '
$ws.Range("D20").Value = '11.47'
$ws.Range("E20").Value = '  +4.46%  '
$ws.Range("D21").Value = '7.75'
$ws.Range("E21").Value = '  +3.13%  '
$ws.Range("D22").Value = '4.13'
$ws.Range("E22").Value = '  +1.25%  '
$ws.Range("D23").Value = '72.04'
$ws.Range("E23").Value = '  +1.32%  '
$ws.Range("D24").Value = '4.42'
$ws.Range("E24").Value = '  +4.43%  '
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("E26").Value = '  +6.53%  '
$ws.Range("D27").Value = '9.58'
$ws.Range("E27").Value = '  +5.78%  '
$ws.Range("E28").Value = '  +5.61%  '
$ws.Range("E29").Value = '  +0.37%  '
$ws.Range("D30").Value = '0.0₃0949'
$ws.Range("E30").Value = '  +5.39%  '
$ws.Range("D31").Value = '526.61'
$ws.Range("E31").Value = '  +3.53%  '
$ws.Range("E32").Value = '  +2.14%  '
$ws.Range("E33").Value = '  +3.19%  '
$ws.Range("E34").Value = '  +2.51%  '
$ws.Range("D36").Value = '165.41'
$ws.Range("E36").Value = '  +1.96%  '
$ws.Range("E37").Value = '  -1.80%  '
$ws.Range("D38").Value = '19.12'
$ws.Range("E38").Value = '  +4.14%  '
$ws.Range("D39").Value = '1.87'
$ws.Range("E39").Value = '  +5.83%  '
$ws.Range("E40").Value = '  +1.42%  '
$ws.Range("E41").Value = '  +2.91%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").Value = '2.60'
$ws.Range("E43").Value = '  +7.77%  '
$ws.Range("D44").Value = '5.01'
$ws.Range("E44").Value = '  +3.67%  '
$ws.Range("D45").Value = '0.329'
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("D46").Value = '40.11'
$ws.Range("E46").Value = '  +2.94%  '
$ws.Range("D47").Value = '153.80'
$ws.Range("E47").Value = '  +2.64%  '
$ws.Range("D48").Value = '3.62'
$ws.Range("E48").Value = '  +1.82%  '
$ws.Range("D49").Value = '0.530'
$ws.Range("E49").Value = '  +2.19%  '
$ws.Range("E50").Value = '  +4.57%  '
$ws.Range("D51").Value = '0.0₆0262'
$ws.Range("E51").Value = '  +0.80%  '
